$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets("ALC")
$ws.Range("H70").Value = 1580.2727
$ws.Range("J70").Value = 1348
$ws.Range("L70").Value = 4044
$ws.Range("N70").Value = -4584

$ws.Range("H73").Value = 1580.2727
$ws.Range("J73").Value = 1348
$ws.Range("L73").Value = 4044
$ws.Range("N73").Value = -5916

$ws.Range("H106").Value = 22379.521
$ws.Range("I106").Value = 6237.25
$ws.Range("K106").Value = 6237.25
$ws.Range("M106").Value = -5606.25

$ws.Range("H137").Value = 16424.193
$ws.Range("I137").Value = 6621.1665
$ws.Range("K137").Value = 19863.4995
$ws.Range("M137").Value = -17313.4995

$ws = $wb.Worksheets("ARM")
$ws.Range("H45").Value = 9338.666999999999
$ws.Range("I45").Value = 13443.667
$ws.Range("K45").Value = 13443.667
$ws.Range("M45").Value = -13066.667

$ws.Range("H61").Value = 9334.241
$ws.Range("I61").Value = 7247.926
$ws.Range("K61").Value = 7247.926
$ws.Range("M61").Value = -7035.926

$ws.Range("H74").Value = 22642.857
$ws.Range("I74").Value = 22642.857
$ws.Range("K74").Value = 22642.857
$ws.Range("M74").Value = -21768.857

$ws.Range("H77").Value = 22642.857
$ws.Range("I77").Value = 22642.857
$ws.Range("K77").Value = 113214.285
$ws.Range("M77").Value = -108846.285

$ws.Range("H110").Value = 813.5
$ws.Range("I110").Value = 882.1667
$ws.Range("J110").Value = 607.5
$ws.Range("K110").Value = 882.1667
$ws.Range("L110").Value = 607.5
$ws.Range("M110").Value = 1162.8333
$ws.Range("N110").Value = -4697.5

$ws.Range("H122").Value = 2311.4783
$ws.Range("I122").Value = 2328.353
$ws.Range("K122").Value = 6985.059
$ws.Range("M122").Value = -4535.059

$ws.Range("H132").Value = 10271.143
$ws.Range("I132").Value = 3649.6667
$ws.Range("K132").Value = 10949.0001
$ws.Range("M132").Value = -8419.000100000001

$ws.Range("H136").Value = 9334.241
$ws.Range("I136").Value = 7247.926
$ws.Range("K136").Value = 21743.778
$ws.Range("M136").Value = -19193.778

$ws = $wb.Worksheets("BSM")
$ws.Range("H107").Value = 1034.7561
$ws.Range("I107").Value = 708.6053000000001
$ws.Range("J107").Value = 5166
$ws.Range("K107").Value = 708.6053000000001
$ws.Range("L107").Value = 5166
$ws.Range("M107").Value = 1211.3947
$ws.Range("N107").Value = -9006

$ws.Range("H134").Value = 6369.7075
$ws.Range("I134").Value = 2979.6667
$ws.Range("J134").Value = 15941.588
$ws.Range("K134").Value = 8939.000100000001
$ws.Range("L134").Value = 47824.764
$ws.Range("M134").Value = -6404.000100000001
$ws.Range("N134").Value = -52894.764

$ws = $wb.Worksheets("CRP")
$ws.Range("H31").Value = 2819.9285
$ws.Range("I31").Value = 1287.3334
$ws.Range("J31").Value = 3237.9092
$ws.Range("K31").Value = 1287.3334
$ws.Range("L31").Value = 3237.9092
$ws.Range("M31").Value = -992.3334
$ws.Range("N31").Value = -3827.9092

$ws.Range("H34").Value = 2819.9285
$ws.Range("I34").Value = 1287.3334
$ws.Range("J34").Value = 3237.9092
$ws.Range("K34").Value = 1287.3334
$ws.Range("L34").Value = 3237.9092
$ws.Range("M34").Value = -1085.3334
$ws.Range("N34").Value = -3641.9092

$ws.Range("H58").Value = 11799.286
$ws.Range("I58").Value = 4898.6665
$ws.Range("J58").Value = 16974.75
$ws.Range("K58").Value = 4898.6665
$ws.Range("L58").Value = 16974.75
$ws.Range("M58").Value = -4695.6665
$ws.Range("N58").Value = -17380.75

$ws.Range("H136").Value = 11799.286
$ws.Range("I136").Value = 4898.6665
$ws.Range("J136").Value = 16974.75
$ws.Range("K136").Value = 14695.9995
$ws.Range("L136").Value = 50924.25
$ws.Range("M136").Value = -12145.9995
$ws.Range("N136").Value = -56024.25

$ws = $wb.Worksheets("GSM")
$ws.Range("H7").Value = 46.666668
$ws.Range("I7").Value = 20
$ws.Range("J7").Value = 100
$ws.Range("K7").Value = 20
$ws.Range("L7").Value = 100
$ws.Range("M7").Value = 92
$ws.Range("N7").Value = -324

$ws.Range("H8").Value = 46.666668
$ws.Range("I8").Value = 20
$ws.Range("J8").Value = 100
$ws.Range("K8").Value = 20
$ws.Range("L8").Value = 100
$ws.Range("M8").Value = 119
$ws.Range("N8").Value = -378

$ws.Range("H80").Value = 2939.6
$ws.Range("I80").Value = 2899.3333
$ws.Range("K80").Value = 2899.3333
$ws.Range("M80").Value = -1901.3333

$ws.Range("H83").Value = 2939.6
$ws.Range("I83").Value = 2899.3333
$ws.Range("K83").Value = 14496.6665
$ws.Range("M83").Value = -9504.666499999999

$ws.Range("H102").Value = 1803.125
$ws.Range("I102").Value = 2024.8182
$ws.Range("J102").Value = 1315.4
$ws.Range("K102").Value = 2024.8182
$ws.Range("L102").Value = 1315.4
$ws.Range("M102").Value = -402.8181999999999
$ws.Range("N102").Value = -4559.4

$ws.Range("H107").Value = 612.2
$ws.Range("J107").Value = 628.1667
$ws.Range("L107").Value = 628.1667
$ws.Range("N107").Value = -4468.1667

$ws.Range("H132").Value = 123798.6
$ws.Range("I132").Value = 109996.5
$ws.Range("J132").Value = 133000
$ws.Range("K132").Value = 329989.5
$ws.Range("L132").Value = 399000
$ws.Range("M132").Value = -327459.5
$ws.Range("N132").Value = -404060

$ws = $wb.Worksheets("LTW")
$ws.Range("H16").Value = 25579.5
$ws.Range("I16").Value = 17127.6
$ws.Range("K16").Value = 17127.6
$ws.Range("M16").Value = -16957.6

$ws.Range("H61").Value = 858.6667
$ws.Range("I61").Value = 858.6667
$ws.Range("K61").Value = 858.6667
$ws.Range("M61").Value = -656.6667

$ws.Range("H113").Value = 858.6667
$ws.Range("I113").Value = 858.6667
$ws.Range("K113").Value = 858.6667
$ws.Range("M113").Value = 1311.3333

$ws.Range("H122").Value = 3318.8
$ws.Range("I122").Value = 2659.6
$ws.Range("J122").Value = 3978
$ws.Range("K122").Value = 7978.799999999999
$ws.Range("L122").Value = 11934
$ws.Range("M122").Value = -5528.799999999999
$ws.Range("N122").Value = -16834

$ws.Range("H136").Value = 6737.054
$ws.Range("I136").Value = 6189.96
$ws.Range("J136").Value = 7876.8335
$ws.Range("K136").Value = 18569.88
$ws.Range("L136").Value = 23630.5005
$ws.Range("M136").Value = -16019.88
$ws.Range("N136").Value = -28730.5005

$ws.Range("H140").Value = 77428.39999999999
$ws.Range("J140").Value = 77428.39999999999
$ws.Range("L140").Value = 77428.39999999999
$ws.Range("N140").Value = -87788.39999999999

$ws = $wb.Worksheets("WVR")
$ws.Range("H136").Value = 7149635
$ws.Range("I136").Value = 12508762
$ws.Range("J136").Value = 4132.6665
$ws.Range("K136").Value = 37526286
$ws.Range("L136").Value = 12397.9995
$ws.Range("M136").Value = -37523736
$ws.Range("N136").Value = -17497.9995
